# Rename each sheet (rerun-with-new-seed naming) and update the
# "Education[T.Unknown]" label to "Education[T.Unknown/Other]" in cell A5
# of every summary sheet.

$wb = $excel.ActiveWorkbook

$newNames = @(
    "summ19178014",
    "summ19414263",
    "summ19678998",
    "summ20003946",
    "summ20308641",
    "summ20615596",
    "summ20881251",
    "summ21147308",
    "summ21428696"
)

for ($i = 1; $i -le $wb.Worksheets.Count; $i++) {
    $ws = $wb.Worksheets.Item($i)
    $ws.Name = $newNames[$i - 1]

    $cell = $ws.Range("A5")
    if ($cell.Value2 -eq "Education[T.Unknown]") {
        $cell.Value2 = "Education[T.Unknown/Other]"
    }
}
